$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025-07-23 09:33:16"
$ws.Range("B2").Value = "add-user"
$ws.Range("C2").Value = "new-organization97"
$ws.Range("D2").Value = "firstteam"
$ws.Range("F2").Value = "Vignesh2122"
$ws.Range("G2").Value = "push"
$ws.Range("I2").Value = "'False"
